$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Closing_Price")

# Update existing values that changed (new Global M2 revision figures)
$ws.Range("B196").Value = 114976911118419.2
$ws.Range("B197").Value = 114672959226172
$ws.Range("B198").Value = 112917914904701.7

# Append new row 199 (2023-06-01) with its Global M2 value,
# copying the date formatting used by the rest of column A
$ws.Range("A198").Copy($ws.Range("A199"))
$ws.Range("A199").Value = 45078
$ws.Range("B199").Value = 112569394359339.9
